# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# "Bad Drivers" table: refreshed critical-minute counts / roaming % for
# the Intel AX211 23.90.0.2 row (and its rollup in the Totals row).
$ws.Range("C3").Value = 2212
$ws.Range("D3").Value = 92.59999999999999
$ws.Range("C4").Value = 2212

# "Good Drivers" table: E12 (driver vintage for AX211 22.150.3.1) was
# previously blank; fill in the date it was observed. A leading
# apostrophe forces Excel to store it as literal text ("2022-08-29")
# instead of auto-parsing it into a date serial number, matching the
# neighboring E13/E14 cells which are already plain text dates. Copying
# E13's format over afterwards clears the "number stored as text"
# quote-prefix flag that typing in the apostrophe sets, so E12 ends up
# back on the same shared style as the rest of column E/D in this block.
$ws.Range("E12").Value = "'2022-08-29"
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
